# Weekly Fruta/Hortaliza update: insert 4 new price rows (date 44588) at the
# top of the "Femacal de La Calera - Sandia" block, pushing the existing
# rows 326-410 down to 330-414.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert four blank rows before the current row 326 (format is inherited
# from the row above, matching native Excel "Insert" behaviour).
$ws.Range("326:329").EntireRow.Insert()

# Common, constant values shared by the four new rows.
$mercadoId = 3
$mercado   = "Femacal de La Calera"
$region    = "Coquimbo"
$codreg    = 5
$catId     = 100112028
$categoria = "Sandia"
$variedad  = "Sin especificar"
$unidadCom = "`$/unidad"
$origen    = "Paine"
$kgOUnid   = 1
$clasif    = "Hortaliza"
$fecha     = 44588

# Row, calidad, volumen, precio minimo, precio maximo, precio promedio
$filas = @(
    @(326, "Extra",   300, 3000, 3000, 3000),
    @(327, "Primera", 350, 2000, 2000, 2000),
    @(328, "Segunda", 560, 1500, 1600, 1554),
    @(329, "Tercera", 380, 1000, 1000, 1000)
)

foreach ($fila in $filas) {
    $r = $fila[0]
    $ws.Range("A$r").Value = $mercadoId
    $ws.Range("B$r").Value = $mercado
    $ws.Range("C$r").Value = $region
    $ws.Range("D$r").Value = $fecha
    $ws.Range("E$r").Value = $codreg
    $ws.Range("F$r").Value = $catId
    $ws.Range("G$r").Value = $categoria
    $ws.Range("H$r").Value = $variedad
    $ws.Range("I$r").Value = $fila[1]
    $ws.Range("J$r").Value = $fila[2]
    $ws.Range("K$r").Value = $fila[3]
    $ws.Range("L$r").Value = $fila[4]
    $ws.Range("M$r").Value = $fila[5]
    $ws.Range("N$r").Value = $unidadCom
    $ws.Range("O$r").Value = $origen
    $ws.Range("P$r").Value = $fila[5]
    $ws.Range("Q$r").Value = $kgOUnid
    $ws.Range("R$r").Value = $clasif
}
